# Update data: 8 October 2021
# Adds the September 2021 unemployment observations to both sheets:
#   - "Canada"   sheet: one new row (Canada, national figure)
#   - "Province" sheet: ten new rows (one per province/territory)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Canada": append row 22
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Canada")

$dateFmt = $ws1.Range("A21").NumberFormat

$ws1.Range("A22").Value = 44440
$ws1.Range("A22").NumberFormat = $dateFmt

$ws1.Range("B22").Value = "Canada"
$ws1.Range("B22").NumberFormat = $dateFmt

$ws1.Range("C22").Formula = "=(D22-E22)/E22*100"

$ws1.Range("D22").Value = 1421.8
$ws1.Range("E22").Value = 1124.4000000000001

# ---------------------------------------------------------------------
# Sheet "Province": append rows 202-211 (one per province)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Province")

$dateFmt2 = $ws2.Range("A201").NumberFormat

$provinces = @(
    @{ Row = 202; Name = "Newfoundland & Labrador"; D = 33.299999999999997;  E = 29.6 },
    @{ Row = 203; Name = "Prince Edward Island";     D = 10.1;               E = 7.3 },
    @{ Row = 204; Name = "Nova Scotia";              D = 40.700000000000003; E = 36.9 },
    @{ Row = 205; Name = "New Brunswick";            D = 37.1;               E = 32.200000000000003 },
    @{ Row = 206; Name = "Quebec";                   D = 261.7;              E = 231.4 },
    @{ Row = 207; Name = "Ontario";                  D = 591.6;              E = 419.8 },
    @{ Row = 208; Name = "Manitoba";                 D = 39.1;               E = 33.700000000000003 },
    @{ Row = 209; Name = "Saskatchewan";              D = 38.4;              E = 33.299999999999997 },
    @{ Row = 210; Name = "Alberta";                  D = 201.2;              E = 164.4 },
    @{ Row = 211; Name = "British Columbia";         D = 168.6;              E = 135.80000000000001 }
)

foreach ($p in $provinces) {
    $r = $p.Row

    $ws2.Cells.Item($r, 1).Value = 44440
    $ws2.Cells.Item($r, 1).NumberFormat = $dateFmt2

    $ws2.Cells.Item($r, 2).Value = $p.Name
    if ($r -eq 202) {
        $ws2.Cells.Item($r, 2).NumberFormat = $dateFmt2
    }

    $ws2.Cells.Item($r, 3).Formula = "=(D$r-E$r)/E$r*100"

    $ws2.Cells.Item($r, 4).Value = $p.D
    $ws2.Cells.Item($r, 5).Value = $p.E
}

# ---------------------------------------------------------------------
# Update view/selection state to match the committed workbook
# ---------------------------------------------------------------------
$ws1.Range("A22").Select()
$ws2.Range("D212").Select()
$ws2.Activate()
